# Apply the commit's data changes to the weather/outfit log sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the old rows 5 and 6 (from the bottom up so row indices stay valid).
$ws.Rows.Item(6).Delete()
$ws.Rows.Item(5).Delete()

# Row 3 now holds what used to be row 5's data (with the comment changed to " None").
# Leading apostrophe keeps the date-looking string stored as text (not an auto-converted date).
$ws.Range("A3").Value = "'01-03-2024"
$ws.Range("B3").Value = 1
$ws.Range("C3").Value = 42.67
$ws.Range("D3").Value = 32.32
$ws.Range("E3").Value = 55
$ws.Range("F3").Value = 8.1
$ws.Range("G3").Value = "scattered clouds"
$ws.Range("H3").Value = ""
$ws.Range("I3").Value = "heavy hoddie"
$ws.Range("J3").Value = "joggers"
$ws.Range("K3").Value = "boots"
$ws.Range("L3").Value = 10
$ws.Range("M3").Value = " None"

# Row 4 becomes a brand-new entry for 01-05-2024.
$ws.Range("A4").Value = "'01-05-2024"
$ws.Range("B4").Value = 1
$ws.Range("C4").Value = 36.54
$ws.Range("D4").Value = 27.82
$ws.Range("E4").Value = 42
$ws.Range("F4").Value = 10.09
$ws.Range("G4").Value = "clear sky"
$ws.Range("H4").Value = "hat"
$ws.Range("I4").Value = "tshirt"
$ws.Range("J4").Value = "joggers"
$ws.Range("K4").Value = "sneakers"
$ws.Range("L4").Value = "'5"
$ws.Range("M4").Value = "cold"
